# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# to reflect the newer snapshot scraped at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 35
$wsExhibition.Range("F4").Value = 217
$wsExhibition.Range("F5").Value = 2645
$wsExhibition.Range("F6").Value = 1881
$wsExhibition.Range("F9").Value = 936

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 35
$wsAll.Range("F4").Value = 217
$wsAll.Range("F5").Value = 2645
$wsAll.Range("F6").Value = 1881
$wsAll.Range("F10").Value = 936
